$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 426; this shifts the old rows 426..490 down to 427..491
# and keeps the dimension/row count consistent with the target (A1:R491).
$ws.Rows.Item(426).Insert()

# Match the date-formatted style used by column D in the surrounding rows.
$ws.Range("D426").NumberFormat = $ws.Range("D427").NumberFormat

# Populate the newly inserted row 426 with its data.
$ws.Range("A426").Value = 10
$ws.Range("B426").Value = "Vega Modelo de Temuco"
$ws.Range("C426").Value = "La Araucanía"
$ws.Range("D426").Value = 44474
$ws.Range("E426").Value = 9
$ws.Range("F426").Value = 100112006
$ws.Range("G426").Value = "Repollo"
$ws.Range("H426").Value = "Crespo record"
$ws.Range("I426").Value = "Primera"
$ws.Range("J426").Value = 1000
$ws.Range("K426").Value = 1000
$ws.Range("L426").Value = 1000
$ws.Range("M426").Value = 1000
$ws.Range("N426").Value = "`$/unidad"
$ws.Range("O426").Value = "Región Metropolitana"
$ws.Range("P426").Value = 1000
$ws.Range("Q426").Value = 1
$ws.Range("R426").Value = "Hortaliza"
